# Normalize "Unidad de Medida" (column X) labels and (re)compute
# "Toneladas Finales" (column Y) from "Cantidad Comercial" (column W).
#
# Units that are converted straight across (already expressed in tonnes):
#   TONELADA       -> TONELADAS        (Y = W)
# Units expressed in kilograms that get normalized to KILOGRAMOS and
# converted to tonnes (Y = W / 1000):
#   KILOGRAMO        -> KILOGRAMOS
#   KILOS NETOS      -> KILOGRAMOS
#   KILOGRAMO BRUTO  -> KILOGRAMOS
#   KG               -> KILOGRAMOS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unitMap = @{
    "TONELADA"        = "TONELADAS"
    "KILOGRAMO"       = "KILOGRAMOS"
    "KILOS NETOS"     = "KILOGRAMOS"
    "KILOGRAMO BRUTO" = "KILOGRAMOS"
    "KG"              = "KILOGRAMOS"
}

for ($r = 2; $r -le 253; $r++) {
    $unitCell = $ws.Cells.Item($r, 24)   # column X
    $unit = $unitCell.Value2

    if ($null -eq $unit) { continue }
    if (-not $unitMap.ContainsKey($unit)) { continue }

    $qtyCell = $ws.Cells.Item($r, 23)    # column W
    $qty = $qtyCell.Value2
    if ($null -eq $qty) { continue }

    if ($unit -eq "TONELADA") {
        $tonnes = $qty
    } else {
        $tonnes = $qty / 1000
    }

    $unitCell.Value = $unitMap[$unit]
    $ws.Cells.Item($r, 25).Value = $tonnes   # column Y
}
